$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet ships protected (legacy password-hash lock on edits). Unprotect so
# the Weight/Percent Change refresh below can land, then restore protection
# with the same allowances afterwards.
$ws.Unprotect()

# Shared string: bump the "as of" model date by one day (2021-05-10 -> 2021-05-11).
$ws.Range("A38").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-11 for illustrative purposes only and are subject to change."

# Refreshed Weight (D) / Percent Change (E) figures for each holding row (2-35, incl. Total row 35).
$ws.Range("D2").Value = 0.03594249038619183
$ws.Range("E2").Value = -0.003488372093023218
$ws.Range("D3").Value = 0.02034975963131538
$ws.Range("E3").Value = -0.00349901553780374
$ws.Range("D4").Value = 0.0190886792352927
$ws.Range("E4").Value = -0.001613553852359839
$ws.Range("D5").Value = 0.03751190314470377
$ws.Range("E5").Value = -0.01544943820224731
$ws.Range("D6").Value = 0.03405448379816789
$ws.Range("E6").Value = -0.0004001600640255454
$ws.Range("D7").Value = 0.01966198601993052
$ws.Range("E7").Value = -0.001352657004830893
$ws.Range("D8").Value = 0.03732117238176648
$ws.Range("E8").Value = -0.00610140963601935
$ws.Range("D9").Value = 0.02030712712643954
$ws.Range("E9").Value = -0.002334979793443992
$ws.Range("D10").Value = 0.02627677204132164
$ws.Range("E10").Value = -0.001060343165606326
$ws.Range("D11").Value = 0.02398065734540318
$ws.Range("E11").Value = -0.01109057301293903
$ws.Range("D12").Value = 0.05762044737852334
$ws.Range("E12").Value = -0.01161971830985919
$ws.Range("D13").Value = 0.02492113984578324
$ws.Range("E13").Value = -0.007683863885839637
$ws.Range("D14").Value = 0.02754537224463526
$ws.Range("E14").Value = -0.007691147639873441
$ws.Range("D15").Value = 0.03377043933925968
$ws.Range("E15").Value = -0.01646505376344098
$ws.Range("D16").Value = 0.01958685886948271
$ws.Range("E16").Value = -0.01279478173607618
$ws.Range("D17").Value = 0.0312236137414085
$ws.Range("E17").Value = -0.005597468970552466
$ws.Range("D18").Value = 0.04178295510306858
$ws.Range("E18").Value = -0.004837595024187902
$ws.Range("D19").Value = 0.1252719650972167
$ws.Range("E19").Value = -0.005976095617529875
$ws.Range("D20").Value = 0.00910801450179537
$ws.Range("E20").Value = -0.009082768325444635
$ws.Range("D21").Value = 0.01551074236285302
$ws.Range("E21").Value = -0.01159122085048003
$ws.Range("D22").Value = 0.01696480480338136
$ws.Range("E22").Value = -0.004111883851972142
$ws.Range("D23").Value = 0.0156764476286081
$ws.Range("E23").Value = -0.01127157449806271
$ws.Range("D24").Value = 0.02148489591994
$ws.Range("E24").Value = -0.01135557132718235
$ws.Range("D25").Value = 0.0124123933078708
$ws.Range("E25").Value = -0.01011758271807495
$ws.Range("D26").Value = 0.04235134797429208
$ws.Range("E26").Value = -0.00802023899235671
$ws.Range("D27").Value = 0.02375983519516581
$ws.Range("E27").Value = 0.0001961553550411388
$ws.Range("D28").Value = 0.04555684060526696
$ws.Range("E28").Value = -0.00710563713879675
$ws.Range("D29").Value = 0.05545427276412791
$ws.Range("E29").Value = -0.006624888093106573
$ws.Range("D30").Value = 0.01299248838526655
$ws.Range("E30").Value = -0.01229773462783168
$ws.Range("D31").Value = 0.02047835421242309
$ws.Range("E31").Value = -0.0007668711656441118
$ws.Range("D32").Value = 0.01361457969173933
$ws.Range("E32").Value = -0.007906976744186056
$ws.Range("D33").Value = 0.04161456491555103
$ws.Range("E33").Value = -0.000514933058702316
$ws.Range("D34").Value = 0.0168025950018076
$ws.Range("E34").Value = -0.01152312010636736
$ws.Range("D35").Value = 1
$ws.Range("E35").Value = -0.00682395966389171

# Restore sheet protection (objects/scenarios locked, column/row formatting still allowed).
$ws.Protect($null, $true, $true, $true, $false, $false, $false, $true, $false, $false, $false, $false, $false, $false, $false, $false)
